$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13-18 shift down to 14-19.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new review entry.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "`nRapaz atencioso e cordial. Obrigado "
$ws.Range("C13").Value = 45926.50978196759
$ws.Range("D13").Value = "ZjZhZGU5MzktZGFiNC00YjdlLTgxN2EtYTkyYWM1ZWE5YjAzOjU3MDE2"

# Avoid a stray custom row height from the embedded newline in B13.
$ws.Rows.Item(13).EntireRow.AutoFit()
